$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.114.59"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.136.59"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D5").Value = "571.25"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "161.82"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -5.81%  "
$ws.Range("D9").Value = "3.149.44"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").Value = "6.59"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "3.686.78"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "64.202.55"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "24.98"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "3.134.97"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "0.0000154"
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "401.65"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "5.24"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "7.13"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "5.86"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "67.74"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").Value = "0.483"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000100"
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("D29").Value = "8.79"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "0.989"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "1.81"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").Value = "21.14"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").Value = "159.22"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").Value = "4.81"
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("D37").Value = "1.11"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").Value = "1.34"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "2.664.06"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").Value = "1.67"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "23.63"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "4.07"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "38.31"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "0.688"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "0.0612"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "5.43"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "288.83"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "0.0255"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").Value = "21.07"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D51").Value = "0.0977"
$ws.Range("E51").Value = "  -0.86%  "
